# Generate Report for Handback
# Updates the localization-status workbook to reflect a completed handback:
#  - Status moves from "Ready for handoff" to "Handed back: in sync with en-US"
#  - Latest Handback DateTime is refreshed
#  - The stale "handback file is not latest" Error Detail is cleared
#  - A couple of columns are widened so the new, longer text is readable

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsOverview = $wb.Worksheets.Item("Overview")

$newStatus = "Handed back: in sync with en-US"

# --- Overview sheet, data row 2 ---
# E2/F2 mirror the per-language Status value, so they move to the new status too
$wsOverview.Range("E2").Value2 = $newStatus
$wsOverview.Range("F2").Value2 = $newStatus

# --- zh-cn sheet, data row 2 ---
$wsZhCn.Range("C2").Value2 = $newStatus
$wsZhCn.Range("K2").Value2 = "2016-08-31 04:51:12"
$wsZhCn.Range("P2").Value2 = ""

# --- de-de sheet, data row 2 ---
$wsDeDe.Range("C2").Value2 = $newStatus
$wsDeDe.Range("K2").Value2 = "2016-08-31 04:51:19"
$wsDeDe.Range("P2").Value2 = ""

# --- widen columns that now hold the longer status / shorter error text ---
# Overview: columns E (zh-cn) and F (de-de) show the Status text
$wsOverview.Columns.Item(5).ColumnWidth = 29.166666666666668
$wsOverview.Columns.Item(6).ColumnWidth = 29.166666666666668

# zh-cn / de-de: column C is Status (wider), column P is Error Detail (narrower now that it is empty)
$wsZhCn.Columns.Item(3).ColumnWidth = 29.166666666666668
$wsZhCn.Columns.Item(16).ColumnWidth = 12.833333333333334

$wsDeDe.Columns.Item(3).ColumnWidth = 29.166666666666668
$wsDeDe.Columns.Item(16).ColumnWidth = 12.833333333333334
